$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = $text
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue 2 4 '="59.355.56"'
Set-TextValue 2 5 '="  -3.86%  "'
Set-TextValue 3 4 '="2.368.52"'
Set-TextValue 3 5 '="  -3.16%  "'
Set-TextValue 4 5 '="  +0.24%  "'
Set-TextValue 5 4 '="556.17"'
Set-TextValue 5 5 '="  -3.70%  "'
Set-TextValue 6 4 '="136.38"'
Set-TextValue 6 5 '="  -3.36%  "'
Set-TextValue 7 5 '="  +0.13%  "'
Set-TextValue 8 4 '="0.530"'
Set-TextValue 8 5 '="  -0.32%  "'
Set-TextValue 9 4 '="2.368.44"'
Set-TextValue 9 5 '="  -2.88%  "'
Set-TextValue 10 4 '="0.104"'
Set-TextValue 10 5 '="  -5.77%  "'
Set-TextValue 11 4 '="0.159"'
Set-TextValue 11 5 '="  -1.37%  "'
Set-TextValue 12 4 '="5.01"'
Set-TextValue 12 5 '="  -3.11%  "'
Set-TextValue 13 4 '="0.334"'
Set-TextValue 13 5 '="  -2.27%  "'
Set-TextValue 14 4 '="25.24"'
Set-TextValue 14 5 '="  -3.38%  "'
Set-TextValue 15 4 '="2.815.34"'
Set-TextValue 15 5 '="  -2.73%  "'
Set-TextValue 16 4 '="0.0000162"'
Set-TextValue 16 5 '="  -5.30%  "'
Set-TextValue 17 4 '="59.435.10"'
Set-TextValue 17 5 '="  -3.67%  "'
Set-TextValue 18 4 '="8.28"'
Set-TextValue 18 5 '="  +13.45%  "'
Set-TextValue 19 4 '="2.362.12"'
Set-TextValue 19 5 '="  -2.44%  "'
Set-TextValue 20 4 '="10.42"'
Set-TextValue 20 5 '="  -1.99%  "'
Set-TextValue 21 4 '="319.22"'
Set-TextValue 21 5 '="  -2.05%  "'
Set-TextValue 22 4 '="3.98"'
Set-TextValue 22 5 '="  -2.18%  "'
Set-TextValue 23 4 '="5.97"'
Set-TextValue 23 5 '="  +0.06%  "'
Set-TextValue 24 5 '="  +0.01%  "'
Set-TextValue 25 4 '="1.76"'
Set-TextValue 25 5 '="  -8.55%  "'
Set-TextValue 26 4 '="63.90"'
Set-TextValue 26 5 '="  -1.80%  "'
Set-TextValue 27 4 '="541.19"'
Set-TextValue 27 5 '="  -7.13%  "'
Set-TextValue 28 4 '="7.91"'
Set-TextValue 28 5 '="  -13.67%  "'
Set-TextValue 29 4 '="2.504.18"'
Set-TextValue 29 5 '="  -2.52%  "'
Set-TextValue 30 4 '="0.0₃0891"'
Set-TextValue 30 5 '="  -3.39%  "'
Set-TextValue 31 4 '="7.78"'
Set-TextValue 31 5 '="  -1.87%  "'
Set-TextValue 32 5 '="  -5.99%  "'
Set-TextValue 33 4 '="1.76"'
Set-TextValue 33 5 '="  -5.63%  "'
Set-TextValue 34 4 '="0.128"'
Set-TextValue 34 5 '="  -4.32%  "'
Set-TextValue 35 5 '="  -0.18%  "'
Set-TextValue 36 4 '="152.31"'
Set-TextValue 36 5 '="  +0.87%  "'
Set-TextValue 37 4 '="1.39"'
Set-TextValue 37 5 '="  -0.08%  "'
Set-TextValue 38 4 '="0.362"'
Set-TextValue 38 5 '="  -2.69%  "'
Set-TextValue 39 4 '="4.44"'
Set-TextValue 39 5 '="  -5.18%  "'
Set-TextValue 40 4 '="18.02"'
Set-TextValue 40 5 '="  -2.00%  "'
Set-TextValue 41 4 '="4.92"'
Set-TextValue 41 5 '="  -4.21%  "'
Set-TextValue 42 5 '="  -0.03%  "'
Set-TextValue 43 4 '="41.02"'
Set-TextValue 43 5 '="  -1.74%  "'
Set-TextValue 44 4 '="1.62"'
Set-TextValue 44 5 '="  -3.76%  "'
Set-TextValue 45 4 '="2.26"'
Set-TextValue 45 5 '="  -4.50%  "'
Set-TextValue 46 4 '="0.0₆0280"'
Set-TextValue 46 5 '="  -5.17%  "'
Set-TextValue 47 4 '="138.89"'
Set-TextValue 47 5 '="  -2.84%  "'
Set-TextValue 48 4 '="3.44"'
Set-TextValue 48 5 '="  -3.45%  "'
Set-TextValue 49 4 '="0.581"'
Set-TextValue 49 5 '="  -2.76%  "'
Set-TextValue 50 4 '="0.0493"'
Set-TextValue 50 5 '="  -3.25%  "'
Set-TextValue 51 4 '="18.67"'
Set-TextValue 51 5 '="  -5.18%  "'

$excel.CutCopyMode = 0

